# Edit script: applies the documented OOXML diff via the Word COM object model.
$d = $word.ActiveDocument

# --- 1) Remove the single comment thread (author POLACK Francois) ---------
# This also strips the commentRangeStart/commentRangeEnd/commentReference
# markers that surround "par" in the paragraph below.
if ($d.Comments.Count -gt 0) {
    for ($i = $d.Comments.Count; $i -ge 1; $i--) {
        $d.Comments.Item($i).Delete()
    }
}

# --- 2) Merge " étant défini " + "par" into a single run " étant défini par"
$null = $d.Content.Find.Execute(" étant défini par", $true, $false, $false, $false, $false, `
    $true, 1, $false, " étant défini par", 2)

# --- 3) Append the new "NB :" paragraph content (text + inline equation) --
# to the last (empty) paragraph of the document, just before the sectPr.
$p = $d.Paragraphs.Last
$rng = $p.Range
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="40D28DE0" w14:textId="77777777" w:rsidR="00E65B25" w:rsidRPr="00E65B25" w:rsidRDefault="00E65B25" w:rsidP="005E5C8E"><w:pPr><w:tabs><w:tab w:val="right" w:pos="8505"/></w:tabs></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">NB : </w:t></w:r><w:r><w:t xml:space="preserve">Dans le programme </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>OptiX</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, la section elliptique est décrite dans le plan XZ. Dans l’espace 3D, la rotation d’angle</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:sym w:font="Symbol" w:char="F06A"/></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve">de X vers Z est une rotation d’angle </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:sym w:font="Symbol" w:char="F06A"/></w:r><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t xml:space="preserve"> et de vecteur </w:t></w:r><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>-</m:t></m:r><m:acc><m:accPr><m:chr m:val="⃗"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:accPr><m:e><m:sSub><m:sSubPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:sSubPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>e</m:t></m:r></m:e><m:sub><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsiaTheme="minorEastAsia" w:hAnsi="Cambria Math"/></w:rPr><m:t>y</m:t></m:r></m:sub></m:sSub></m:e></m:acc></m:oMath><w:r><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$rng.InsertXML($xml)
